$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title cell (A1, merged A1:E1)
$ws.Range("A1").Value = "PURCHASE STATUS REPORT"

# Add new "Status" header in F4, to the right of the existing header row
$ws.Range("F4").Value = "Status"

# Update the selection to match the new active cell (F4)
$ws.Range("F4").Select()
